$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New Harris Interactive poll (11/17) — three rows appended at the bottom
# of the data table (rows 108-110), following the same column layout as
# the existing rows.
#
# Columns: A id, B year, C week, D month, E day, F firm, G collectmode,
# H excluded, I n, J unsure, K c_poutou, L c_arthaud, M c_melenchon,
# N c_roussel, O c_jadot, P c_hidalgo, Q c_macron, R c_pecresse,
# S c_barnier, T c_bertrand, U c_lassalle, V c_daignant, W c_lepen,
# X c_zemmour, Y c_asselineau, Z c_poisson, AA c_philippot, AB c_lagarde

$rows = @(
    @{ Row=108; A=32; B=2021; C=12; D=11; E=14; F="harris"; G="online"; H="included"; I=2027; J=1;   K=1; L=10; M=2; N=2; O=8; P=4; Q=23;      T=14;     U="T_0.5"; V=1; W=16; X=17; Y="T_0.5"; AA=1 },
    @{ Row=109; A=32; B=2021; C=12; D=11; E=14; F="harris"; G="online"; H="included"; I=2028; J=1;   K=1; L=10; M=2; N=2; O=9; P=4; Q=24; R=10;        U="T_0.5"; V=2; W=16; X=17; Y="T_0.5"; AA=2 },
    @{ Row=110; A=32; B=2021; C=12; D=11; E=14; F="harris"; G="online"; H="included"; I=2029; J=1;   K=1; L=10; M=2; N=2; O=9; P=4; Q=24;      S=10;     U="T_0.5"; V=2; W=16; X=18; Y="T_0.5"; AA=1 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    if ($r.ContainsKey("R")) { $ws.Cells.Item($row, 18).Value = $r.R }
    if ($r.ContainsKey("S")) { $ws.Cells.Item($row, 19).Value = $r.S }
    if ($r.ContainsKey("T")) { $ws.Cells.Item($row, 20).Value = $r.T }
    $ws.Cells.Item($row, 21).Value = $r.U
    $ws.Cells.Item($row, 22).Value = $r.V
    $ws.Cells.Item($row, 23).Value = $r.W
    $ws.Cells.Item($row, 24).Value = $r.X
    # column Y (c_asselineau) mirrors the "T_0.5" styling used elsewhere
    # (explicit black font, style index 1)
    $yCell = $ws.Cells.Item($row, 25)
    $yCell.Value = $r.Y
    $yCell.Font.Color = 0
    $ws.Cells.Item($row, 27).Value = $r.AA
}

# Selection state, matching the final view after data entry
$ws.Range("AB110").Select()
